# Template + Version update
# Updates the "script_IQSlidedeck.R | Date: ..." footer timestamps that are
# re-generated on each run of the R reporting script.

$p = $ppt.ActivePresentation

$updates = @(
    @{ Slide = 3;  Footer = "Footer Placeholder 4"; Old = "2021-01-11 16:15:43"; New = "2021-01-11 16:38:09" },
    @{ Slide = 5;  Footer = "Footer Placeholder 5"; Old = "2021-01-11 16:15:43"; New = "2021-01-11 16:38:10" },
    @{ Slide = 6;  Footer = "Footer Placeholder 4"; Old = "2021-01-11 16:15:47"; New = "2021-01-11 16:38:13" },
    @{ Slide = 7;  Footer = "Footer Placeholder 4"; Old = "2021-01-11 16:15:48"; New = "2021-01-11 16:38:14" },
    @{ Slide = 8;  Footer = "Footer Placeholder 4"; Old = "2021-01-11 16:15:54"; New = "2021-01-11 16:38:21" },
    @{ Slide = 10; Footer = "Footer Placeholder 5"; Old = "2021-01-11 16:15:44"; New = "2021-01-11 16:38:10" },
    @{ Slide = 11; Footer = "Footer Placeholder 5"; Old = "2021-01-11 16:15:45"; New = "2021-01-11 16:38:12" },
    @{ Slide = 12; Footer = "Footer Placeholder 5"; Old = "2021-01-11 16:15:46"; New = "2021-01-11 16:38:13" }
)

foreach ($u in $updates) {
    $s = $p.Slides.Item($u.Slide)
    $shp = $s.Shapes.Item($u.Footer)
    $tr = $shp.TextFrame.TextRange
    $tr.Text = $tr.Text.Replace($u.Old, $u.New)
}
